$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$base = "/home/masoud/Documents/four-polar/fourPolar-io/target/test-classes/fr/fresnel/fourPolar/io/imageSet/acquisition/sample/finders/excel/SampleImageSetByExcelFileFinder/OneCamera"

$ws.Range("A5").Value = "$base/Img1_C1.tif"
$ws.Range("A6").Value = "$base/Img2_C1.tif"
$ws.Range("A7").Value = "$base/Img3_C1.tif"
$ws.Range("A8").Value = "$base/Img4_C1.tif"

$ws.Range("A5:A8").Select()
